# Update line-loading percentage results for the 380 kV case (Case_0_235).
# Only the C, D, E, F, G, H, K, M columns for rows 2-25 change; all other
# cells (A, B, I, J, L, N, O and row 1 headers) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4.91766526950903
$ws.Range("D2").Value = 4.981210237019123
$ws.Range("E2").Value = 10.48592781916091
$ws.Range("F2").Value = 27.83069089757192
$ws.Range("G2").Value = 35.48657044786039
$ws.Range("H2").Value = 15.47493534562134
$ws.Range("K2").Value = 13.9202465219519
$ws.Range("M2").Value = 15.64267658331061

$ws.Range("C3").Value = 4.749957266423106
$ws.Range("D3").Value = 5.001471611759303
$ws.Range("E3").Value = 10.4052058710073
$ws.Range("F3").Value = 27.71133405364528
$ws.Range("G3").Value = 35.2015840284863
$ws.Range("H3").Value = 15.50583129137171
$ws.Range("K3").Value = 13.39447410288986
$ws.Range("M3").Value = 15.347152538532

$ws.Range("C4").Value = 4.645541737644154
$ws.Range("D4").Value = 5.01445172421634
$ws.Range("E4").Value = 10.35945688099465
$ws.Range("F4").Value = 27.64905176804063
$ws.Range("G4").Value = 35.04295172079882
$ws.Range("H4").Value = 15.52924793218144
$ws.Range("K4").Value = 13.06341005362698
$ws.Range("M4").Value = 15.16689637838861

$ws.Range("C5").Value = 4.602705024714634
$ws.Range("D5").Value = 5.019877539076165
$ws.Range("E5").Value = 10.34178744332683
$ws.Range("F5").Value = 27.62644938766646
$ws.Range("G5").Value = 34.98247488134486
$ws.Range("H5").Value = 15.53990236918454
$ws.Range("K5").Value = 12.9266404698908
$ws.Range("M5").Value = 15.09384459039849

$ws.Range("C6").Value = 4.595577051027535
$ws.Range("D6").Value = 5.020786743458962
$ws.Range("E6").Value = 10.33891262976668
$ws.Range("F6").Value = 27.62286430446181
$ws.Range("G6").Value = 34.97268584087993
$ws.Range("H6").Value = 15.54173850351284
$ws.Range("K6").Value = 12.90382440756847
$ws.Range("M6").Value = 15.08174182828547

$ws.Range("C7").Value = 4.644965082115554
$ws.Range("D7").Value = 5.014524345904699
$ws.Range("E7").Value = 10.3592146251435
$ws.Range("F7").Value = 27.64873568477811
$ws.Range("G7").Value = 35.0421191713755
$ws.Range("H7").Value = 15.52938712782382
$ws.Range("K7").Value = 13.06157275761199
$ws.Range("M7").Value = 15.16590940363769

$ws.Range("C8").Value = 4.860182722848325
$ws.Range("D8").Value = 4.988084759980606
$ws.Range("E8").Value = 10.4573125119656
$ws.Range("F8").Value = 27.78726122191241
$ws.Range("G8").Value = 35.38495032769522
$ws.Range("H8").Value = 15.48466185658723
$ws.Range("K8").Value = 13.74079103139953
$ws.Range("M8").Value = 15.54059306191339

$ws.Range("C9").Value = 5.267738042052621
$ws.Range("D9").Value = 4.940488939000141
$ws.Range("E9").Value = 10.67920840921199
$ws.Range("F9").Value = 28.14548059683442
$ws.Range("G9").Value = 36.1839983430492
$ws.Range("H9").Value = 15.43249172665903
$ws.Range("K9").Value = 14.99910741710948
$ws.Range("M9").Value = 16.28037867574993

$ws.Range("C10").Value = 5.5546926282475
$ws.Range("D10").Value = 4.908072414760643
$ws.Range("E10").Value = 10.85910260573028
$ws.Range("F10").Value = 28.46012086738486
$ws.Range("G10").Value = 36.84364840893138
$ws.Range("H10").Value = 15.41614870333955
$ws.Range("K10").Value = 15.86909682558553
$ws.Range("M10").Value = 16.82114963420917

$ws.Range("C11").Value = 5.681877156237681
$ws.Range("D11").Value = 4.893871097228104
$ws.Range("E11").Value = 10.94432587644875
$ws.Range("F11").Value = 28.61405297996716
$ws.Range("G11").Value = 37.15832541869133
$ws.Range("H11").Value = 15.41354719884832
$ws.Range("K11").Value = 16.25140000652258
$ws.Range("M11").Value = 17.06539615097843

$ws.Range("C12").Value = 5.729510222182395
$ws.Range("D12").Value = 4.888571194436159
$ws.Range("E12").Value = 10.97705952591981
$ws.Range("F12").Value = 28.67385774059125
$ws.Range("G12").Value = 37.27947594673374
$ws.Range("H12").Value = 15.41326105425813
$ws.Range("K12").Value = 16.39411895405599
$ws.Range("M12").Value = 17.15754221056843

$ws.Range("C13").Value = 5.71927580344021
$ws.Range("D13").Value = 4.889709170673998
$ws.Range("E13").Value = 10.96998964878853
$ws.Range("F13").Value = 28.66091104365029
$ws.Range("G13").Value = 37.25329735459663
$ws.Range("H13").Value = 15.41329153491908
$ws.Range("K13").Value = 16.36347476816032
$ws.Range("M13").Value = 17.13771366793912

$ws.Range("C14").Value = 5.685806807495501
$ws.Range("D14").Value = 4.893433514410085
$ws.Range("E14").Value = 10.9470097972781
$ws.Range("F14").Value = 28.61894304857505
$ws.Range("G14").Value = 37.16825332965002
$ws.Range("H14").Value = 15.41350962779007
$ws.Range("K14").Value = 16.26318329147532
$ws.Range("M14").Value = 17.07298448699154

$ws.Range("C15").Value = 5.665235897227967
$ws.Range("D15").Value = 4.895724901587974
$ws.Range("E15").Value = 10.93299328311882
$ws.Range("F15").Value = 28.59343238059135
$ws.Range("G15").Value = 37.11641709587658
$ws.Range("H15").Value = 15.41373435449048
$ws.Range("K15").Value = 16.20148160561455
$ws.Range("M15").Value = 17.03328842200239

$ws.Range("C16").Value = 5.54630922783429
$ws.Range("D16").Value = 4.909011425334687
$ws.Range("E16").Value = 10.85359906971116
$ws.Range("F16").Value = 28.45027512648164
$ws.Range("G16").Value = 36.82336800814757
$ws.Range("H16").Value = 15.41641632538449
$ws.Range("K16").Value = 15.84383134215997
$ws.Range("M16").Value = 16.80514424160112

$ws.Range("C17").Value = 5.472457812792284
$ws.Range("D17").Value = 4.917301493425963
$ws.Range("E17").Value = 10.80574317917611
$ws.Range("F17").Value = 28.36519254879121
$ws.Range("G17").Value = 36.64725449129129
$ws.Range("H17").Value = 15.41930250245697
$ws.Range("K17").Value = 15.62088998440122
$ws.Range("M17").Value = 16.66466874167432

$ws.Range("C18").Value = 5.429667316516912
$ws.Range("D18").Value = 4.922121063438701
$ws.Range("E18").Value = 10.77853830612198
$ws.Range("F18").Value = 28.31727451797794
$ws.Range("G18").Value = 36.54733976466785
$ws.Range("H18").Value = 15.42141725629032
$ws.Range("K18").Value = 15.49139900580424
$ws.Range("M18").Value = 16.58371121555008

$ws.Range("C19").Value = 5.415126929429929
$ws.Range("D19").Value = 4.923761722351537
$ws.Range("E19").Value = 10.76938302883797
$ws.Range("F19").Value = 28.30122647365365
$ws.Range("G19").Value = 36.51375082361067
$ws.Range("H19").Value = 15.42221123990687
$ws.Range("K19").Value = 15.44734278536472
$ws.Range("M19").Value = 16.5562758568201

$ws.Range("C20").Value = 5.480352184949506
$ws.Range("D20").Value = 4.916413691974392
$ws.Range("E20").Value = 10.81080452603761
$ws.Range("F20").Value = 28.37414452032138
$ws.Range("G20").Value = 36.66585991647314
$ws.Range("H20").Value = 15.41894817358552
$ws.Range("K20").Value = 15.64475386840001
$ws.Range("M20").Value = 16.67963981721618

$ws.Range("C21").Value = 5.695652165269989
$ws.Range("D21").Value = 4.892337476540636
$ws.Range("E21").Value = 10.95374722449255
$ws.Range("F21").Value = 28.63122930738495
$ws.Range("G21").Value = 37.19317974024659
$ws.Range("H21").Value = 15.41342656905441
$ws.Range("K21").Value = 16.2926978360271
$ws.Range("M21").Value = 17.09200708926662

$ws.Range("C22").Value = 5.833261995476451
$ws.Range("D22").Value = 4.877055668516155
$ws.Range("E22").Value = 11.04984548787285
$ws.Range("F22").Value = 28.80805305536169
$ws.Range("G22").Value = 37.54935349740733
$ws.Range("H22").Value = 15.41389359106364
$ws.Range("K22").Value = 16.70416907378566
$ws.Range("M22").Value = 17.3594636580484

$ws.Range("C23").Value = 5.760114897520202
$ws.Range("D23").Value = 4.88517055329105
$ws.Range("E23").Value = 10.99831990277388
$ws.Range("F23").Value = 28.71288724592813
$ws.Range("G23").Value = 37.35823869862158
$ws.Range("H23").Value = 15.41327022344375
$ws.Range("K23").Value = 16.48569114790169
$ws.Range("M23").Value = 17.21693393928913

$ws.Range("C24").Value = 5.476784170865815
$ws.Range("D24").Value = 4.916814900296074
$ws.Range("E24").Value = 10.80851532946083
$ws.Range("F24").Value = 28.37009422524315
$ws.Range("G24").Value = 36.65744423617966
$ws.Range("H24").Value = 15.41910694724643
$ws.Range("K24").Value = 15.63396910598044
$ws.Range("M24").Value = 16.67287199826606

$ws.Range("C25").Value = 5.159434407370208
$ws.Range("D25").Value = 4.952913842013427
$ws.Range("E25").Value = 10.6161283301143
$ws.Range("F25").Value = 28.03942208005284
$ws.Range("G25").Value = 35.95471248876668
$ws.Range("H25").Value = 15.44276557783166
$ws.Range("K25").Value = 14.66765119747791
$ws.Range("M25").Value = 16.08032061652337
